# Revise the metadata files
#
# 1) preparation_medium sheet: add two new medium options (NBF / PLP),
#    rename "PFA" -> "PFA (Paraformaldehyde)", add a 22nd row
#    ("PAXgene tissue kit (PXT)").
# 2) storage_medium sheet: add the same two new medium options, rename
#    "PFA" -> "PFA (Paraformaldehyde)" and
#    "Paraffin embedded (FFPE)" -> "FFPE (Paraffin embedded)".
# 3) .metadata sheet: bump the pav:createdOn timestamp.
# 4) Sample Section data-validation range for preparation_medium grows
#    from $A$1:$A$21 to $A$1:$A$22.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# preparation_medium (22 rows after the edit)
# ---------------------------------------------------------------------
$wsPrepMedium = $wb.Worksheets.Item("preparation_medium")

$prepMediumData = New-Object 'object[,]' 22,2

$prepMediumData[0,0]  = "PBS"
$prepMediumData[0,1]  = "http://purl.obolibrary.org/obo/OBI_0100046"
$prepMediumData[1,0]  = "NBF (Neutral Buffered Formalin)"
$prepMediumData[1,1]  = "http://purl.obolibrary.org/obo/OBIB_0000213"
$prepMediumData[2,0]  = "PLP (Periodate-lysine-paraformaldehyde)"
$prepMediumData[2,1]  = "http://purl.bioontology.org/ontology/MESH/C046311"
$prepMediumData[3,0]  = "Ethanol"
$prepMediumData[3,1]  = "http://purl.obolibrary.org/obo/CHEBI_16236"
$prepMediumData[4,0]  = "Allprotect tissue reagent (ALL)"
$prepMediumData[4,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118"
$prepMediumData[5,0]  = "CLARITY hydrogel"
$prepMediumData[5,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000134"
$prepMediumData[6,0]  = "Inflated (Agarose)"
$prepMediumData[6,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000106"
$prepMediumData[7,0]  = "MACS tissue storage solution"
$prepMediumData[7,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105"
$prepMediumData[8,0]  = "Fresh frozen CMC"
$prepMediumData[8,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000130"
$prepMediumData[9,0]  = "Inflated (OCT)"
$prepMediumData[9,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000123"
$prepMediumData[10,0] = "Fresh frozen gelatin"
$prepMediumData[10,1] = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000198"
$prepMediumData[11,0] = "PFA (Paraformaldehyde)"
$prepMediumData[11,1] = "http://purl.obolibrary.org/obo/CHEBI_61538"
$prepMediumData[12,0] = "Fixed frozen OCT (Formalin, sucrose protected)"
$prepMediumData[12,1] = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000116"
$prepMediumData[13,0] = "Unknown"
$prepMediumData[13,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"
$prepMediumData[14,0] = "Fresh frozen OCT"
$prepMediumData[14,1] = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000126"
$prepMediumData[15,0] = "RNAlater"
$prepMediumData[15,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"
$prepMediumData[16,0] = "Fixed frozen OCT (Cytofix/Cytoperm)"
$prepMediumData[16,1] = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000149"
$prepMediumData[17,0] = "None"
$prepMediumData[17,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132"
$prepMediumData[18,0] = "Bouin's"
$prepMediumData[18,1] = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000140"
$prepMediumData[19,0] = "Fixed frozen OCT (PFA, sucrose protected)"
$prepMediumData[19,1] = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000147"
$prepMediumData[20,0] = "Methanol"
$prepMediumData[20,1] = "http://purl.obolibrary.org/obo/CHEBI_17790"
$prepMediumData[21,0] = "PAXgene tissue kit (PXT)"
$prepMediumData[21,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113"

$wsPrepMedium.Range("A1:B22").Value = $prepMediumData

# ---------------------------------------------------------------------
# storage_medium (18 rows, same count, reshuffled + renamed content)
# ---------------------------------------------------------------------
$wsStorageMedium = $wb.Worksheets.Item("storage_medium")

$storageMediumData = New-Object 'object[,]' 18,2

$storageMediumData[0,0]  = "PBS"
$storageMediumData[0,1]  = "http://purl.obolibrary.org/obo/OBI_0100046"
$storageMediumData[1,0]  = "OCT"
$storageMediumData[1,1]  = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63523"
$storageMediumData[2,0]  = "NBF (Neutral Buffered Formalin)"
$storageMediumData[2,1]  = "http://purl.obolibrary.org/obo/OBIB_0000213"
$storageMediumData[3,0]  = "PLP (Periodate-lysine-paraformaldehyde)"
$storageMediumData[3,1]  = "http://purl.bioontology.org/ontology/MESH/C046311"
$storageMediumData[4,0]  = "Allprotect tissue reagent (ALL)"
$storageMediumData[4,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000118"
$storageMediumData[5,0]  = "DMSO (no serum)"
$storageMediumData[5,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000115"
$storageMediumData[6,0]  = "MACS tissue storage solution"
$storageMediumData[6,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000105"
$storageMediumData[7,0]  = "PFA (Paraformaldehyde)"
$storageMediumData[7,1]  = "http://purl.obolibrary.org/obo/CHEBI_61538"
$storageMediumData[8,0]  = "Tris-EDTA"
$storageMediumData[8,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000135"
$storageMediumData[9,0]  = "DMSO (with serum)"
$storageMediumData[9,1]  = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125"
$storageMediumData[10,0] = "Unknown"
$storageMediumData[10,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"
$storageMediumData[11,0] = "Gelatin"
$storageMediumData[11,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C65802"
$storageMediumData[12,0] = "RNAlater"
$storageMediumData[12,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"
$storageMediumData[13,0] = "FFPE (Paraffin embedded)"
$storageMediumData[13,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C143028"
$storageMediumData[14,0] = "CMC"
$storageMediumData[14,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C83594"
$storageMediumData[15,0] = "None"
$storageMediumData[15,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C41132"
$storageMediumData[16,0] = "Methanol"
$storageMediumData[16,1] = "http://purl.obolibrary.org/obo/CHEBI_17790"
$storageMediumData[17,0] = "PAXgene tissue kit (PXT)"
$storageMediumData[17,1] = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C185113"

$wsStorageMedium.Range("A1:B18").Value = $storageMediumData

# ---------------------------------------------------------------------
# .metadata: bump pav:createdOn
# ---------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item(".metadata")
$wsMetadata.Range("C2").Value = "2023-08-04T07:35:19-07:00"

# ---------------------------------------------------------------------
# Sample Section: grow the preparation_medium validation list range
# ---------------------------------------------------------------------
$wsSample = $wb.Worksheets.Item("Sample Section")
$prepValidation = $wsSample.Range("D2:D1001").Validation
$prepValidation.Modify(3, 1, 1, "'preparation_medium'!`$A`$1:`$A`$22")
$prepValidation.IgnoreBlank = $true
$prepValidation.ShowError = $true
$prepValidation.ErrorTitle = "Validation Error"
$prepValidation.ErrorMessage = ""
